# Edit script: insert two new price rows (date 2023-01-06 / serial 44924)
# for "Vega Monumental Concepción - Cilantro" at the top of the data block,
# pushing all existing rows (189 onward) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 189 (existing data shifts down to make room)
$ws.Rows.Item(189).Insert()
$ws.Rows.Item(189).Insert()

# New row 189: "Primera" quality entry for the new reporting date
$ws.Range("A189").Value = 11
$ws.Range("B189").Value = 'Vega Monumental Concepción'
$ws.Range("C189").Value = 'Bíobío'
$ws.Range("D189").Value = 44924
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = 100112040
$ws.Range("G189").Value = 'Cilantro'
$ws.Range("H189").Value = 'Sin especificar'
$ws.Range("I189").Value = 'Primera'
$ws.Range("J189").Value = 200
$ws.Range("K189").Value = 700
$ws.Range("L189").Value = 800
$ws.Range("M189").Value = 750
$ws.Range("N189").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O189").Value = 'Región de Ñuble'
$ws.Range("P189").Value = 750
$ws.Range("Q189").Value = 1
$ws.Range("R189").Value = 'Hortaliza'

# New row 190: "Segunda" quality entry for the new reporting date
$ws.Range("A190").Value = 11
$ws.Range("B190").Value = 'Vega Monumental Concepción'
$ws.Range("C190").Value = 'Bíobío'
$ws.Range("D190").Value = 44924
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = 100112040
$ws.Range("G190").Value = 'Cilantro'
$ws.Range("H190").Value = 'Sin especificar'
$ws.Range("I190").Value = 'Segunda'
$ws.Range("J190").Value = 100
$ws.Range("K190").Value = 600
$ws.Range("L190").Value = 600
$ws.Range("M190").Value = 600
$ws.Range("N190").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O190").Value = 'Región de Ñuble'
$ws.Range("P190").Value = 600
$ws.Range("Q190").Value = 1
$ws.Range("R190").Value = 'Hortaliza'

# Apply the same date-number format used by the other "Fecha" (column D) cells
$ws.Range("D189").NumberFormat = $ws.Range("D191").NumberFormat
$ws.Range("D190").NumberFormat = $ws.Range("D191").NumberFormat
